# Fruta / hortaliza, semanal
# Rotate the "variable" fields (Fecha, Calidad, Volumen, Precio min/max/prom,
# Unidad de comercialización, Precio $/Kg, Kg/unidad) across rows 2-6 so that
# each row picks up the data that used to belong three rows further down
# (wrapping around within the 2-6 block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that move, keyed by row.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")
$snapshot = @{}
foreach ($r in 2..6) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# New row r gets the snapshot of row (source) below, wrapping 2..6.
$mapping = @{ 2 = 5; 3 = 6; 4 = 2; 5 = 3; 6 = 4 }

foreach ($r in 2..6) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $snapshot[$src][$c]
    }
}
